$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 6424.75
$ws.Cells.Item(19, 9).Value = 6499.25
$ws.Cells.Item(19, 11).Value = 6499.25
$ws.Cells.Item(19, 13).Value = -6324.25
$ws.Cells.Item(38, 8).Value = 31.222221
$ws.Cells.Item(38, 9).Value = 31.222221
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 93.666663
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = 278.333337
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 1313
$ws.Cells.Item(40, 10).Value = 1313
$ws.Cells.Item(40, 12).Value = 1313
$ws.Cells.Item(40, 14).Value = -1663
$ws.Cells.Item(62, 8).Value = 5693.25
$ws.Cells.Item(62, 9).Value = 4942.5
$ws.Cells.Item(62, 11).Value = 4942.5
$ws.Cells.Item(62, 13).Value = -4318.5
$ws.Cells.Item(65, 8).Value = 5693.25
$ws.Cells.Item(65, 9).Value = 4942.5
$ws.Cells.Item(65, 11).Value = 24712.5
$ws.Cells.Item(65, 13).Value = -21592.5
$ws.Cells.Item(98, 8).Value = 892.6
$ws.Cells.Item(98, 9).Value = 1159.25
$ws.Cells.Item(98, 10).Value = 587.8570999999999
$ws.Cells.Item(98, 11).Value = 1159.25
$ws.Cells.Item(98, 12).Value = 587.8570999999999
$ws.Cells.Item(98, 13).Value = 338.75
$ws.Cells.Item(98, 14).Value = -3583.8571
$ws.Cells.Item(103, 8).Value = 1862.8823
$ws.Cells.Item(103, 9).Value = 599.8333
$ws.Cells.Item(103, 11).Value = 1799.4999
$ws.Cells.Item(103, 13).Value = -1213.4999
$ws.Cells.Item(122, 8).Value = 892.6
$ws.Cells.Item(122, 9).Value = 1159.25
$ws.Cells.Item(122, 10).Value = 587.8570999999999
$ws.Cells.Item(122, 11).Value = 3477.75
$ws.Cells.Item(122, 12).Value = 1763.5713
$ws.Cells.Item(122, 13).Value = -1027.75
$ws.Cells.Item(122, 14).Value = -6663.5713
$ws.Cells.Item(132, 8).Value = 3930.1155
$ws.Cells.Item(132, 9).Value = 3133.125
$ws.Cells.Item(132, 10).Value = 5205.3
$ws.Cells.Item(132, 11).Value = 9399.375
$ws.Cells.Item(132, 12).Value = 15615.9
$ws.Cells.Item(132, 13).Value = -6869.375
$ws.Cells.Item(132, 14).Value = -20675.9
$ws.Cells.Item(138, 8).Value = 2684.5
$ws.Cells.Item(138, 9).Value = 2746.7273
$ws.Cells.Item(138, 11).Value = 8240.1819
$ws.Cells.Item(138, 13).Value = -3100.1819
$ws.Cells.Item(141, 8).Value = 1126.7646
$ws.Cells.Item(141, 9).Value = 1150.3334
$ws.Cells.Item(141, 11).Value = 3451.0002
$ws.Cells.Item(141, 13).Value = 1728.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 1334.9375
$ws.Cells.Item(74, 9).Value = 1334.9375
$ws.Cells.Item(74, 11).Value = 1334.9375
$ws.Cells.Item(74, 13).Value = -460.9375
$ws.Cells.Item(77, 8).Value = 1334.9375
$ws.Cells.Item(77, 9).Value = 1334.9375
$ws.Cells.Item(77, 11).Value = 6674.6875
$ws.Cells.Item(77, 13).Value = -2306.6875
$ws.Cells.Item(97, 8).Value = 776.88
$ws.Cells.Item(97, 9).Value = 530.5454999999999
$ws.Cells.Item(97, 10).Value = 2583.3333
$ws.Cells.Item(97, 11).Value = 530.5454999999999
$ws.Cells.Item(97, 12).Value = 2583.3333
$ws.Cells.Item(97, 13).Value = -34.54549999999995
$ws.Cells.Item(97, 14).Value = -3575.3333
$ws.Cells.Item(122, 8).Value = 3210.077
$ws.Cells.Item(122, 9).Value = 3304.2
$ws.Cells.Item(122, 10).Value = 2896.3333
$ws.Cells.Item(122, 11).Value = 9912.599999999999
$ws.Cells.Item(122, 12).Value = 8688.999899999999
$ws.Cells.Item(122, 13).Value = -7462.599999999999
$ws.Cells.Item(122, 14).Value = -13588.9999
$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3478.25
$ws.Cells.Item(20, 10).Value = 5665.6665
$ws.Cells.Item(20, 12).Value = 5665.6665
$ws.Cells.Item(20, 14).Value = -6159.6665
$ws.Cells.Item(86, 8).Value = 2761.8
$ws.Cells.Item(86, 9).Value = 2116.75
$ws.Cells.Item(86, 11).Value = 2116.75
$ws.Cells.Item(86, 13).Value = -993.75
$ws.Cells.Item(89, 8).Value = 2761.8
$ws.Cells.Item(89, 9).Value = 2116.75
$ws.Cells.Item(89, 11).Value = 10583.75
$ws.Cells.Item(89, 13).Value = -4967.75
$ws.Cells.Item(94, 8).Value = 2448
$ws.Cells.Item(94, 9).Value = 1960.6666
$ws.Cells.Item(94, 10).Value = 3701.1428
$ws.Cells.Item(94, 11).Value = 1960.6666
$ws.Cells.Item(94, 12).Value = 3701.1428
$ws.Cells.Item(94, 13).Value = -1509.6666
$ws.Cells.Item(94, 14).Value = -4603.1428
$ws.Cells.Item(105, 8).Value = 3866.8
$ws.Cells.Item(105, 9).Value = 2611.3333
$ws.Cells.Item(105, 11).Value = 2611.3333
$ws.Cells.Item(105, 13).Value = -864.3332999999998
$ws.Cells.Item(134, 8).Value = 1340.2106
$ws.Cells.Item(134, 9).Value = 1262.5883
$ws.Cells.Item(134, 10).Value = 2000
$ws.Cells.Item(134, 11).Value = 3787.7649
$ws.Cells.Item(134, 12).Value = 6000
$ws.Cells.Item(134, 13).Value = -1252.7649
$ws.Cells.Item(134, 14).Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 26742.875
$ws.Cells.Item(12, 9).Value = 3485.75
$ws.Cells.Item(12, 11).Value = 3485.75
$ws.Cells.Item(12, 13).Value = -3315.75
$ws.Cells.Item(31, 8).Value = 1890.5
$ws.Cells.Item(31, 9).Value = 1501
$ws.Cells.Item(31, 10).Value = 2474.75
$ws.Cells.Item(31, 11).Value = 1501
$ws.Cells.Item(31, 12).Value = 2474.75
$ws.Cells.Item(31, 13).Value = -1206
$ws.Cells.Item(31, 14).Value = -3064.75
$ws.Cells.Item(34, 8).Value = 1890.5
$ws.Cells.Item(34, 9).Value = 1501
$ws.Cells.Item(34, 10).Value = 2474.75
$ws.Cells.Item(34, 11).Value = 1501
$ws.Cells.Item(34, 12).Value = 2474.75
$ws.Cells.Item(34, 13).Value = -1299
$ws.Cells.Item(34, 14).Value = -2878.75
$ws.Cells.Item(58, 8).Value = 1423.7693
$ws.Cells.Item(58, 9).Value = 1229.3636
$ws.Cells.Item(58, 10).Value = 2493
$ws.Cells.Item(58, 11).Value = 1229.3636
$ws.Cells.Item(58, 12).Value = 2493
$ws.Cells.Item(58, 13).Value = -1026.3636
$ws.Cells.Item(58, 14).Value = -2899
$ws.Cells.Item(86, 8).Value = 13943891
$ws.Cells.Item(89, 8).Value = 13943891
$ws.Cells.Item(132, 8).Value = 2739.8696
$ws.Cells.Item(132, 9).Value = 2858
$ws.Cells.Item(132, 10).Value = 1499.5
$ws.Cells.Item(132, 11).Value = 8574
$ws.Cells.Item(132, 12).Value = 4498.5
$ws.Cells.Item(132, 13).Value = -6044
$ws.Cells.Item(132, 14).Value = -9558.5
$ws.Cells.Item(134, 8).Value = 2585.0557
$ws.Cells.Item(134, 9).Value = 2242.0667
$ws.Cells.Item(134, 11).Value = 6726.2001
$ws.Cells.Item(134, 13).Value = -4191.2001
$ws.Cells.Item(136, 8).Value = 1423.7693
$ws.Cells.Item(136, 9).Value = 1229.3636
$ws.Cells.Item(136, 10).Value = 2493
$ws.Cells.Item(136, 11).Value = 3688.0908
$ws.Cells.Item(136, 12).Value = 7479
$ws.Cells.Item(136, 13).Value = -1138.0908
$ws.Cells.Item(136, 14).Value = -12579

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 651.4
$ws.Cells.Item(5, 9).Value = 696
$ws.Cells.Item(5, 10).Value = 584.5
$ws.Cells.Item(5, 11).Value = 2088
$ws.Cells.Item(5, 12).Value = 1753.5
$ws.Cells.Item(5, 13).Value = -1976
$ws.Cells.Item(5, 14).Value = -1977.5
$ws.Cells.Item(40, 8).Value = 98.25
$ws.Cells.Item(40, 9).Value = 108.6
$ws.Cells.Item(40, 10).Value = 81
$ws.Cells.Item(40, 11).Value = 434.4
$ws.Cells.Item(40, 12).Value = 324
$ws.Cells.Item(40, 13).Value = -365.4
$ws.Cells.Item(40, 14).Value = -462
$ws.Cells.Item(49, 8).Value = 4666.6665
$ws.Cells.Item(49, 9).Value = 5000
$ws.Cells.Item(49, 10).Value = 4600
$ws.Cells.Item(49, 11).Value = 15000
$ws.Cells.Item(49, 12).Value = 13800
$ws.Cells.Item(49, 13).Value = -14844
$ws.Cells.Item(49, 14).Value = -14112
$ws.Cells.Item(122, 8).Value = 591.0714
$ws.Cells.Item(122, 9).Value = 447.85715
$ws.Cells.Item(122, 10).Value = 734.2857
$ws.Cells.Item(122, 11).Value = 4030.71435
$ws.Cells.Item(122, 12).Value = 6608.571300000001
$ws.Cells.Item(122, 13).Value = -1580.71435
$ws.Cells.Item(122, 14).Value = -11508.5713
$ws.Cells.Item(135, 8).Value = 651.4
$ws.Cells.Item(135, 9).Value = 696
$ws.Cells.Item(135, 10).Value = 584.5
$ws.Cells.Item(135, 11).Value = 6264
$ws.Cells.Item(135, 12).Value = 5260.5
$ws.Cells.Item(135, 13).Value = -3729
$ws.Cells.Item(135, 14).Value = -10330.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4690.2856
$ws.Cells.Item(70, 9).Value = 4639.6665
$ws.Cells.Item(70, 11).Value = 4639.6665
$ws.Cells.Item(70, 13).Value = -4369.6665
$ws.Cells.Item(73, 8).Value = 4690.2856
$ws.Cells.Item(73, 9).Value = 4639.6665
$ws.Cells.Item(73, 11).Value = 4639.6665
$ws.Cells.Item(73, 13).Value = -3703.6665
$ws.Cells.Item(97, 8).Value = 1118.7778
$ws.Cells.Item(97, 10).Value = 1557.25
$ws.Cells.Item(97, 12).Value = 1557.25
$ws.Cells.Item(97, 14).Value = -2549.25
$ws.Cells.Item(102, 8).Value = 1202.1428
$ws.Cells.Item(102, 9).Value = 1013.75
$ws.Cells.Item(102, 10).Value = 1453.3334
$ws.Cells.Item(102, 11).Value = 1013.75
$ws.Cells.Item(102, 12).Value = 1453.3334
$ws.Cells.Item(102, 13).Value = 608.25
$ws.Cells.Item(102, 14).Value = -4697.3334
$ws.Cells.Item(132, 8).Value = 2865.2
$ws.Cells.Item(132, 9).Value = 2647.0667
$ws.Cells.Item(132, 10).Value = 3519.6
$ws.Cells.Item(132, 11).Value = 7941.2001
$ws.Cells.Item(132, 12).Value = 10558.8
$ws.Cells.Item(132, 13).Value = -5411.2001
$ws.Cells.Item(132, 14).Value = -15618.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(16, 8).Value = 475.1
$ws.Cells.Item(16, 9).Value = 475.1
$ws.Cells.Item(16, 11).Value = 475.1
$ws.Cells.Item(16, 13).Value = -305.1
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 14).ClearContents()
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 5573.5713
$ws.Cells.Item(40, 9).Value = 5763.4
$ws.Cells.Item(40, 11).Value = 5763.4
$ws.Cells.Item(40, 13).Value = -5627.4
$ws.Cells.Item(46, 8).Value = 1582.4706
$ws.Cells.Item(46, 10).Value = 2007.4286
$ws.Cells.Item(46, 12).Value = 2007.4286
$ws.Cells.Item(46, 14).Value = -2383.4286
$ws.Cells.Item(132, 8).Value = 2912.16
$ws.Cells.Item(132, 9).Value = 2865.25
$ws.Cells.Item(132, 11).Value = 8595.75
$ws.Cells.Item(132, 13).Value = -6065.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 4979262.5
$ws.Cells.Item(100, 9).Value = 9957498
$ws.Cells.Item(100, 11).Value = 19914996
$ws.Cells.Item(100, 13).Value = -19914455
$ws.Cells.Item(122, 8).Value = 4276.6
$ws.Cells.Item(122, 9).Value = 4063.6667
$ws.Cells.Item(122, 11).Value = 12191.0001
$ws.Cells.Item(122, 13).Value = -9741.000100000001
$ws.Cells.Item(132, 8).Value = 2511.1428
$ws.Cells.Item(132, 9).Value = 2393.9375
$ws.Cells.Item(132, 11).Value = 7181.8125
$ws.Cells.Item(132, 13).Value = -4651.8125
